$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2021" data column (O) is being added to the right of the existing
# "2020" column (N). For every row, clone column N's cell formatting (style,
# borders, number format) into column O via copy/paste-special, then fill in
# the 2021 figures. Row 3 (the thin separator row above the header) keeps an
# empty, but formatted, O3 cell - exactly like its N3 counterpart.

$xlPasteFormats = -4122

$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial($xlPasteFormats)

$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial($xlPasteFormats)
$ws.Range("O4").Value = 2021

$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial($xlPasteFormats)
$ws.Range("O5").Value = 70.636215334420882

$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial($xlPasteFormats)
$ws.Range("O6").Value = 107.1

$ws.Range("N7").Copy()
$ws.Range("O7").PasteSpecial($xlPasteFormats)
$ws.Range("O7").Value = 55.452054794520542

$ws.Range("N8").Copy()
$ws.Range("O8").PasteSpecial($xlPasteFormats)
$ws.Range("O8").Value = 84.375

$ws.Range("N9").Copy()
$ws.Range("O9").PasteSpecial($xlPasteFormats)
$ws.Range("O9").Value = 120.48192771084337

$ws.Range("N10").Copy()
$ws.Range("O10").PasteSpecial($xlPasteFormats)
$ws.Range("O10").Value = 109.53346855983774

$ws.Range("N11").Copy()
$ws.Range("O11").PasteSpecial($xlPasteFormats)
$ws.Range("O11").Value = 147.7690288713911

$ws.Range("N12").Copy()
$ws.Range("O12").PasteSpecial($xlPasteFormats)
$ws.Range("O12").Value = 25.545675020210183

$ws.Range("N13").Copy()
$ws.Range("O13").PasteSpecial($xlPasteFormats)
$ws.Range("O13").Value = 82.457854874175425

$ws.Range("N14").Copy()
$ws.Range("O14").PasteSpecial($xlPasteFormats)
$ws.Range("O14").Value = 15.384615384615385

$excel.CutCopyMode = 0
